# Atualização de bases das ligas, do dia: 06-04-2024 às 15:39
# Australia ALeague - refresh of fixture rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 73 and 74 got their fixtures swapped (id column A stays fixed per
#    row position; every other populated column B..AC is exchanged).
# ---------------------------------------------------------------------------
$swapCols73 = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
foreach ($col in $swapCols73) {
    $addrTop = $col + "73"
    $addrBot = $col + "74"
    $valTop = $ws.Range($addrTop).Value()
    $valBot = $ws.Range($addrBot).Value()
    $ws.Range($addrTop).Value = $valBot
    $ws.Range($addrBot).Value = $valTop
}

# ---------------------------------------------------------------------------
# 2) Rows 124 and 125 got their fixtures swapped the same way.
# ---------------------------------------------------------------------------
$swapCols124 = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
foreach ($col in $swapCols124) {
    $addrTop = $col + "124"
    $addrBot = $col + "125"
    $valTop = $ws.Range($addrTop).Value()
    $valBot = $ws.Range($addrBot).Value()
    $ws.Range($addrTop).Value = $valBot
    $ws.Range($addrBot).Value = $valTop
}

# ---------------------------------------------------------------------------
# 3) Row 141 is refreshed with the fixture previously staged in row 144,
#    except the closing odds R/S which came in with new values.
# ---------------------------------------------------------------------------
$copyCols141 = @("B","C","D","E","F","G","K","L","M","N","O","P","Q","T","U","V","W","X","Y","Z","AA")
foreach ($col in $copyCols141) {
    $ws.Range($col + "141").Value = $ws.Range($col + "144").Value()
}
$ws.Range("R141").Value = 1.95
$ws.Range("S141").Value = 1.95

# ---------------------------------------------------------------------------
# 4) Row 142 is refreshed with the fixture previously staged in row 145,
#    except the closing odds R/S/U/V which came in with new values.
# ---------------------------------------------------------------------------
$copyCols142 = @("B","C","D","E","F","G","K","L","M","N","O","P","Q","T","W","X","Y","Z","AA")
foreach ($col in $copyCols142) {
    $ws.Range($col + "142").Value = $ws.Range($col + "145").Value()
}
$ws.Range("R142").Value = 2.06
$ws.Range("S142").Value = 1.84
$ws.Range("U142").Value = 1.925
$ws.Range("V142").Value = 1.925

# ---------------------------------------------------------------------------
# 5) Row 143 becomes a brand new fixture (id/date/teams/odds all new).
# ---------------------------------------------------------------------------
$ws.Range("B143").Value = 7127408
$ws.Range("E143").Value = 45395.10416666666
$ws.Range("F143").Value = "Western United FC"
$ws.Range("G143").Value = "Central Coast Mariners"
$ws.Range("K143").Value = 3.5
$ws.Range("L143").Value = 3.8
$ws.Range("M143").Value = 1.909
$ws.Range("N143").Value = 3.4
$ws.Range("O143").Value = 3.8
$ws.Range("P143").Value = 1.95
$ws.Range("Q143").Value = 0.5
$ws.Range("R143").Value = 1.88
$ws.Range("S143").Value = 2.02
$ws.Range("T143").Value = 3
$ws.Range("U143").Value = 1.975
$ws.Range("V143").Value = 1.875
$ws.Range("W143").Value = 0
$ws.Range("X143").Value = 0
$ws.Range("Y143").Value = 0
$ws.Range("Z143").Value = 0
$ws.Range("AA143").Value = 0

# ---------------------------------------------------------------------------
# 6) Rows 144 and 145 (now fully absorbed into 141-143) are removed, shrinking
#    the sheet from 145 to 143 rows.
# ---------------------------------------------------------------------------
$ws.Range("144:145").EntireRow.Delete()
